$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7921580076217651
$ws.Range("B1").Value = 1.05002498626709
$ws.Range("C1").Value = 1.467407703399658
$ws.Range("D1").Value = 3.115997791290283
$ws.Range("E1").Value = 2.13815450668335
